# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update DAMSLTag (column I) and DialogAct (column J) values for the rows identified
# in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 3;   DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 18;  DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 24;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 25;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 31;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 43;  DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 52;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 59;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 76;  DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 77;  DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 78;  DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 79;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 97;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 101; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 102; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 108; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 120; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 122; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 125; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.DAMSLTag
    $ws.Range("J" + $u.Row).Value = $u.DialogAct
}
